$p = $ppt.ActivePresentation

$oldDate = "26/09/2018"
$newDate = "23/10/2018"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# 1. Update the cached date placeholder text ("26/09/2018" -> "23/10/2018")
#    on every slide master and every slide layout in the deck.
for ($d = 1; $d -le $p.Designs.Count; $d++) {
    $design = $p.Designs.Item($d)
    $master = $design.SlideMaster

    # The master's own date placeholder.
    Update-DatePlaceholder $master.Shapes

    if ($d -eq 1) {
        # Direct navigation works fine for the first design.
        for ($l = 1; $l -le $master.CustomLayouts.Count; $l++) {
            $layout = $master.CustomLayouts.Item($l)
            Update-DatePlaceholder $layout.Shapes
        }
    } else {
        # Direct Master.CustomLayouts navigation for secondary designs does
        # not resolve to the correct underlying layout part in this host,
        # so bind each layout via a throwaway slide (AddSlide binds the
        # CustomLayout reference correctly) and discard the slide again.
        $layoutCount = $master.CustomLayouts.Count
        for ($l = 1; $l -le $layoutCount; $l++) {
            $srcLayout = $master.CustomLayouts.Item($l)
            $tmpSlide = $p.Slides.AddSlide($p.Slides.Count + 1, $srcLayout)
            Update-DatePlaceholder $tmpSlide.CustomLayout.Shapes
            $tmpSlide.Delete()
        }
    }
}

# 2. Bold + underline the "Agile & Waterfall methodologies used" bullet
#    on slide 9.
$slide9 = $p.Slides.Item(9)
$shape = $slide9.Shapes.Item("Content Placeholder 2")
$tf = $shape.TextFrame
for ($para = 1; $para -le $tf.TextRange.Paragraphs().Count; $para++) {
    $run = $tf.TextRange.Paragraphs($para, 1)
    if ($run.Text.TrimEnd("`r`n") -eq "Agile & Waterfall methodologies used") {
        $run.Font.Bold = $true
        $run.Font.Underline = $true
    }
}
